# "went through Programming, found some minor mistakes" -
# the worksheet models an 8-bit binary -> hex converter; cells F5, G5, F7
# and G7 had stray "1" bits checked that shouldn't have been, which threw
# off the downstream DEC2HEX sums. Clear them back to blank/0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()

# leave the cursor parked on F7, the last cell touched
$ws.Range("F7").Select()
